# Auto-generated from diff analysis
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 3
$ws.Range("H3").Value = 29800
$ws.Range("J3").Value = 29800
$ws.Range("L3").Value = 29800
$ws.Range("N3").Value = -30028
# Row 102
$ws.Range("H102").Value = 29800
$ws.Range("J102").Value = 29800
$ws.Range("L102").Value = 29800
$ws.Range("N102").Value = -36290
# Row 112
$ws.Range("H112").Value = 3368074.8
$ws.Range("J112").Value = 3368074.8
$ws.Range("L112").Value = 10104224.4
$ws.Range("N112").Value = -10106440.4
# Row 125
$ws.Range("H125").Value = 687.375
$ws.Range("I125").Value = 474.75
$ws.Range("J125").Value = 900
$ws.Range("K125").Value = 4272.75
$ws.Range("L125").Value = 8100
$ws.Range("M125").Value = -1812.75
$ws.Range("N125").Value = -13020
# Row 135
$ws.Range("H135").Value = 13161637
$ws.Range("I135").Value = 691.8125
$ws.Range("J135").Value = 83353340
$ws.Range("K135").Value = 6226.3125
$ws.Range("L135").Value = 750180060
$ws.Range("M135").Value = -3691.3125
$ws.Range("N135").Value = -750185130
# Row 137
$ws.Range("H137").Value = 1551.3182
$ws.Range("I137").Value = 1227.4166
$ws.Range("J137").Value = 1940
$ws.Range("K137").Value = 3682.2498
$ws.Range("L137").Value = 5820
$ws.Range("M137").Value = -1132.2498
$ws.Range("N137").Value = -10920

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 2148.75
$ws.Range("I2").Value = 1915
$ws.Range("J2").Value = 2850
$ws.Range("K2").Value = 1915
$ws.Range("L2").Value = 2850
$ws.Range("M2").Value = -1802
$ws.Range("N2").Value = -3076
# Row 45
$ws.Range("H45").Value = 2643.0286
$ws.Range("I45").Value = 2395.84
$ws.Range("J45").Value = 3261
$ws.Range("K45").Value = 2395.84
$ws.Range("L45").Value = 3261
$ws.Range("M45").Value = -2018.84
$ws.Range("N45").Value = -4015
# Row 61
$ws.Range("H61").Value = 347076.78
$ws.Range("I61").Value = 410044.06
$ws.Range("K61").Value = 410044.06
$ws.Range("M61").Value = -409832.06
# Row 74
$ws.Range("H74").Value = 26317494
$ws.Range("I74").Value = 28573210
$ws.Range("K74").Value = 28573210
$ws.Range("M74").Value = -28572336
# Row 77
$ws.Range("H77").Value = 26317494
$ws.Range("I77").Value = 28573210
$ws.Range("K77").Value = 142866050
$ws.Range("M77").Value = -142861682
# Row 110
$ws.Range("H110").Value = 979.75
$ws.Range("I110").Value = 962.125
$ws.Range("J110").Value = 1015
$ws.Range("K110").Value = 962.125
$ws.Range("L110").Value = 1015
$ws.Range("M110").Value = 1082.875
$ws.Range("N110").Value = -5105
# Row 116
$ws.Range("H116").Value = 2148.75
$ws.Range("I116").Value = 1915
$ws.Range("J116").Value = 2850
$ws.Range("K116").Value = 1915
$ws.Range("L116").Value = 2850
$ws.Range("M116").Value = 379
$ws.Range("N116").Value = -7438
# Row 122
$ws.Range("H122").Value = 2364.5908
$ws.Range("I122").Value = 1824.0588
$ws.Range("K122").Value = 5472.1764
$ws.Range("M122").Value = -3022.1764
# Row 136
$ws.Range("H136").Value = 347076.78
$ws.Range("I136").Value = 410044.06
$ws.Range("K136").Value = 1230132.18
$ws.Range("M136").Value = -1227582.18

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 2148.75
$ws.Range("I3").Value = 1915
$ws.Range("J3").Value = 2850
$ws.Range("K3").Value = 1915
$ws.Range("L3").Value = 2850
$ws.Range("M3").Value = -1801
$ws.Range("N3").Value = -3078
# Row 61
$ws.Range("H61").Value = 0
$ws.Range("J61").Value = 0
$ws.Range("L61").Value = 0
$ws.Range("N61").ClearContents()
# Row 106
$ws.Range("H106").Value = 17834.2
$ws.Range("J106").Value = 17834.2
$ws.Range("L106").Value = 17834.2
$ws.Range("N106").Value = -20358.2
# Row 134
$ws.Range("H134").Value = 4631.8213
$ws.Range("I134").Value = 5131.409
$ws.Range("J134").Value = 2800
$ws.Range("K134").Value = 15394.227
$ws.Range("L134").Value = 8400
$ws.Range("M134").Value = -12859.227
$ws.Range("N134").Value = -13470

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2314.6614
$ws.Range("I31").Value = 1368.0212
$ws.Range("J31").Value = 5280.8
$ws.Range("K31").Value = 1368.0212
$ws.Range("L31").Value = 5280.8
$ws.Range("M31").Value = -1073.0212
$ws.Range("N31").Value = -5870.8
# Row 34
$ws.Range("H34").Value = 2314.6614
$ws.Range("I34").Value = 1368.0212
$ws.Range("J34").Value = 5280.8
$ws.Range("K34").Value = 1368.0212
$ws.Range("L34").Value = 5280.8
$ws.Range("M34").Value = -1166.0212
$ws.Range("N34").Value = -5684.8
# Row 43
$ws.Range("H43").Value = 39500
$ws.Range("J43").Value = 39500
$ws.Range("L43").Value = 39500
$ws.Range("N43").Value = -39868
# Row 58
$ws.Range("H58").Value = 16550.75
$ws.Range("I58").Value = 981.85187
$ws.Range("J58").Value = 100622.8
$ws.Range("K58").Value = 981.85187
$ws.Range("L58").Value = 100622.8
$ws.Range("M58").Value = -778.85187
$ws.Range("N58").Value = -101028.8
# Row 95
$ws.Range("H95").Value = 32000
$ws.Range("J95").Value = 32000
$ws.Range("L95").Value = 32000
$ws.Range("N95").Value = -37492
# Row 99
$ws.Range("H99").Value = 25003750
$ws.Range("I99").Value = 3500
$ws.Range("J99").Value = 55559612
$ws.Range("K99").Value = 3500
$ws.Range("L99").Value = 55559612
$ws.Range("M99").Value = -2002
$ws.Range("N99").Value = -55562608
# Row 101
$ws.Range("H101").Value = 39500
$ws.Range("J101").Value = 39500
$ws.Range("L101").Value = 39500
$ws.Range("N101").Value = -45990
# Row 126
$ws.Range("H126").Value = 25003750
$ws.Range("I126").Value = 3500
$ws.Range("J126").Value = 55559612
$ws.Range("K126").Value = 10500
$ws.Range("L126").Value = 166678836
$ws.Range("M126").Value = -8030
$ws.Range("N126").Value = -166683776
# Row 132
$ws.Range("H132").Value = 1703.8163
$ws.Range("I132").Value = 1340.7174
$ws.Range("K132").Value = 4022.1522
$ws.Range("M132").Value = -1492.1522
# Row 134
$ws.Range("H134").Value = 794.2759
$ws.Range("J134").Value = 1122.9166
$ws.Range("L134").Value = 3368.7498
$ws.Range("N134").Value = -8438.7498
# Row 136
$ws.Range("H136").Value = 16550.75
$ws.Range("I136").Value = 981.85187
$ws.Range("J136").Value = 100622.8
$ws.Range("K136").Value = 2945.55561
$ws.Range("L136").Value = 301868.4
$ws.Range("M136").Value = -395.5556099999999
$ws.Range("N136").Value = -306968.4

$ws = $wb.Worksheets.Item("CUL")
# Row 131
$ws.Range("H131").Value = 697.4400000000001
$ws.Range("I131").Value = 500
$ws.Range("J131").Value = 714.6087
$ws.Range("K131").Value = 1500
$ws.Range("L131").Value = 2143.8261
$ws.Range("M131").Value = 3540
$ws.Range("N131").Value = -12223.8261
# Row 140
$ws.Range("H140").Value = 1879.4166
$ws.Range("I140").Value = 1405.7
$ws.Range("J140").Value = 4248
$ws.Range("K140").Value = 4217.1
$ws.Range("L140").Value = 12744
$ws.Range("M140").Value = 962.8999999999996
$ws.Range("N140").Value = -23104

$ws = $wb.Worksheets.Item("GSM")
# Row 92
$ws.Range("H92").Value = 13312.625
$ws.Range("I92").Value = 11000
$ws.Range("K92").Value = 11000
$ws.Range("M92").Value = -9128
# Row 101
$ws.Range("H101").Value = 32499.75
$ws.Range("J101").Value = 32499.75
$ws.Range("L101").Value = 32499.75
$ws.Range("N101").Value = -38989.75
# Row 122
$ws.Range("H122").Value = 55556740
$ws.Range("I122").Value = 19608616
$ws.Range("J122").Value = 142859330
$ws.Range("K122").Value = 58825848
$ws.Range("L122").Value = 428577990
$ws.Range("M122").Value = -58823398
$ws.Range("N122").Value = -428582890

$ws = $wb.Worksheets.Item("LTW")
# Row 61
$ws.Range("H61").Value = 6381.727
$ws.Range("I61").Value = 2899.75
$ws.Range("J61").Value = 8371.429
$ws.Range("K61").Value = 2899.75
$ws.Range("L61").Value = 8371.429
$ws.Range("M61").Value = -2697.75
$ws.Range("N61").Value = -8775.429
# Row 93
$ws.Range("H93").Value = 748.6818
$ws.Range("I93").Value = 685.5294
$ws.Range("K93").Value = 685.5294
$ws.Range("M93").Value = 562.4706
# Row 113
$ws.Range("H113").Value = 6381.727
$ws.Range("I113").Value = 2899.75
$ws.Range("J113").Value = 8371.429
$ws.Range("K113").Value = 2899.75
$ws.Range("L113").Value = 8371.429
$ws.Range("M113").Value = -729.75
$ws.Range("N113").Value = -12711.429
# Row 132
$ws.Range("H132").Value = 549054.3
$ws.Range("J132").Value = 800
$ws.Range("L132").Value = 2400
$ws.Range("N132").Value = -7460
# Row 136
$ws.Range("H136").Value = 1039.6285
$ws.Range("I136").Value = 978.6896400000001
$ws.Range("J136").Value = 1334.1666
$ws.Range("K136").Value = 2936.06892
$ws.Range("L136").Value = 4002.4998
$ws.Range("M136").Value = -386.0689200000002
$ws.Range("N136").Value = -9102.4998

$ws = $wb.Worksheets.Item("WVR")
# Row 122
$ws.Range("H122").Value = 1316.5652
$ws.Range("I122").Value = 1335.5
$ws.Range("J122").Value = 900
$ws.Range("K122").Value = 4006.5
$ws.Range("L122").Value = 2700
$ws.Range("M122").Value = -1556.5
$ws.Range("N122").Value = -7600
# Row 132
$ws.Range("H132").Value = 471.69492
$ws.Range("I132").Value = 490.57693
$ws.Range("K132").Value = 1471.73079
$ws.Range("M132").Value = 1058.26921
# Row 136
$ws.Range("H136").Value = 17206300
$ws.Range("I136").Value = 22939954
$ws.Range("K136").Value = 68819862
$ws.Range("M136").Value = -68817312

